# Auto-generated edit script: update cached values in Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 160.25
$ws.Range("I8").Value = 160.25
$ws.Range("K8").Value = 480.75
$ws.Range("M8").Value = -341.75
$ws.Range("H33").Value = 404.30304
$ws.Range("I33").Value = 98.22727
$ws.Range("J33").Value = 1016.4545
$ws.Range("K33").Value = 98.22727
$ws.Range("L33").Value = 1016.4545
$ws.Range("M33").Value = 130.77273
$ws.Range("N33").Value = -1474.4545
$ws.Range("H40").Value = 1495.9166
$ws.Range("J40").Value = 1050
$ws.Range("L40").Value = 1050
$ws.Range("N40").Value = -1400
$ws.Range("H52").Value = 95322.16
$ws.Range("J52").Value = 53716.168
$ws.Range("L52").Value = 161148.504
$ws.Range("N52").Value = -161468.504
$ws.Range("H64").Value = 4900
$ws.Range("H67").Value = 4900
$ws.Range("H70").Value = 3610
$ws.Range("I70").Value = 3567.1428
$ws.Range("K70").Value = 10701.4284
$ws.Range("M70").Value = -10431.4284
$ws.Range("H73").Value = 3610
$ws.Range("I73").Value = 3567.1428
$ws.Range("K73").Value = 10701.4284
$ws.Range("M73").Value = -9765.428400000001
$ws.Range("H93").Value = 24601
$ws.Range("J93").Value = 24601
$ws.Range("L93").Value = 24601
$ws.Range("N93").Value = -29593
$ws.Range("H138").Value = 3336874.5
$ws.Range("I138").Value = 479045.8
$ws.Range("J138").Value = 5212324.5
$ws.Range("K138").Value = 1437137.4
$ws.Range("L138").Value = 15636973.5
$ws.Range("M138").Value = -1431997.4
$ws.Range("N138").Value = -15647253.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4784798.5
$ws.Range("I5").Value = 8772078
$ws.Range("J5").Value = 62.8
$ws.Range("K5").Value = 8772078
$ws.Range("L5").Value = 62.8
$ws.Range("M5").Value = -8771966
$ws.Range("N5").Value = -286.8
$ws.Range("H63").Value = 2633.3333
$ws.Range("I63").Value = 2600
$ws.Range("K63").Value = 2600
$ws.Range("M63").Value = -1914
$ws.Range("H66").Value = 2633.3333
$ws.Range("I66").Value = 2600
$ws.Range("K66").Value = 13000
$ws.Range("M66").Value = -9568
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36622
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -113112

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4784798.5
$ws.Range("I4").Value = 8772078
$ws.Range("J4").Value = 62.8
$ws.Range("K4").Value = 8772078
$ws.Range("L4").Value = 62.8
$ws.Range("M4").Value = -8771963
$ws.Range("N4").Value = -292.8
$ws.Range("H94").Value = 784.1111
$ws.Range("I94").Value = 738.1429000000001
$ws.Range("J94").Value = 945
$ws.Range("K94").Value = 738.1429000000001
$ws.Range("L94").Value = 945
$ws.Range("M94").Value = -287.1429000000001
$ws.Range("N94").Value = -1847

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 142857940
$ws.Range("I22").Value = 200000670
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 200000670
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = -200000320
$ws.Range("N22").Value = -1800
$ws.Range("H62").Value = 2997
$ws.Range("I62").Value = 2999.75
$ws.Range("J62").Value = 2991.5
$ws.Range("K62").Value = 2999.75
$ws.Range("L62").Value = 2991.5
$ws.Range("M62").Value = -2375.75
$ws.Range("N62").Value = -4239.5
$ws.Range("H65").Value = 2997
$ws.Range("I65").Value = 2999.75
$ws.Range("J65").Value = 2991.5
$ws.Range("K65").Value = 14998.75
$ws.Range("L65").Value = 14957.5
$ws.Range("M65").Value = -11878.75
$ws.Range("N65").Value = -21197.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1102.6305
$ws.Range("I131").Value = 715
$ws.Range("J131").Value = 1120.25
$ws.Range("K131").Value = 2145
$ws.Range("L131").Value = 3360.75
$ws.Range("M131").Value = 2895
$ws.Range("N131").Value = -13440.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50749.953
$ws.Range("I70").Value = 88208.336
$ws.Range("J70").Value = 5799.9
$ws.Range("K70").Value = 88208.336
$ws.Range("L70").Value = 5799.9
$ws.Range("M70").Value = -87938.336
$ws.Range("N70").Value = -6339.9
$ws.Range("H73").Value = 50749.953
$ws.Range("I73").Value = 88208.336
$ws.Range("J73").Value = 5799.9
$ws.Range("K73").Value = 88208.336
$ws.Range("L73").Value = 5799.9
$ws.Range("M73").Value = -87272.336
$ws.Range("N73").Value = -7671.9
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -33744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 944
$ws.Range("I22").Value = 712.625
$ws.Range("J22").Value = 1175.375
$ws.Range("K22").Value = 712.625
$ws.Range("L22").Value = 1175.375
$ws.Range("M22").Value = -417.625
$ws.Range("N22").Value = -1765.375
$ws.Range("H27").Value = 944
$ws.Range("I27").Value = 712.625
$ws.Range("J27").Value = 1175.375
$ws.Range("K27").Value = 712.625
$ws.Range("L27").Value = 1175.375
$ws.Range("M27").Value = -605.625
$ws.Range("N27").Value = -1389.375
$ws.Range("H68").Value = 1585.4546
$ws.Range("I68").Value = 1493.3334
$ws.Range("K68").Value = 1493.3334
$ws.Range("M68").Value = -744.3334
$ws.Range("H71").Value = 1585.4546
$ws.Range("I71").Value = 1493.3334
$ws.Range("K71").Value = 7466.666999999999
$ws.Range("M71").Value = -3722.666999999999
$ws.Range("H93").Value = 291.14285
$ws.Range("I93").Value = 291.14285
$ws.Range("K93").Value = 291.14285
$ws.Range("M93").Value = 956.85715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166669980
$ws.Range("J62").Value = 5334.3335
$ws.Range("L62").Value = 5334.3335
$ws.Range("N62").Value = -6582.3335
$ws.Range("H65").Value = 166669980
$ws.Range("J65").Value = 5334.3335
$ws.Range("L65").Value = 26671.6675
$ws.Range("N65").Value = -32911.6675
$ws.Range("H81").Value = 2960.7693
$ws.Range("I81").Value = 5750.5
$ws.Range("J81").Value = 2453.5454
$ws.Range("K81").Value = 11501
$ws.Range("L81").Value = 4907.0908
$ws.Range("M81").Value = -10440
$ws.Range("N81").Value = -7029.0908
$ws.Range("H84").Value = 2960.7693
$ws.Range("I84").Value = 5750.5
$ws.Range("J84").Value = 2453.5454
$ws.Range("K84").Value = 57505
$ws.Range("L84").Value = 24535.454
$ws.Range("M84").Value = -52201
$ws.Range("N84").Value = -35143.454
$ws.Range("H123").Value = 47357.05
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 48876.89
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 48876.89
$ws.Range("M123").Value = -15100
$ws.Range("N123").Value = -58676.89
